# Update LR-pair TPM data (commit: "update scripts wuth new tpm")
# The dataset shrinks from 6 data rows (sender x target combos ECs/FAPs/MuSCs x ECs/FAPs)
# down to 3 data rows (ECs/FAPs/MuSCs -> FAPs), with freshly computed TPM-derived metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> FAPs (sending cluster D2 changes ECs -> FAPs) ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("M2").Value = 0.001809666666666667
$ws.Range("N2").Value = 0.005429
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1156579362958889
$ws.Range("R2").Value = 1.040921426663
$ws.Range("S2").Value = 0.4067926910433548
$ws.Range("T2").Value = 0.4067926910433549

# --- Row 3: sending cluster ECs -> FAPs, target cluster stays FAPs ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 57.4434
$ws.Range("H3").Value = 172.3302
$ws.Range("I3").Value = 0.3656254573230189
$ws.Range("J3").Value = 0.365625457323019
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.1039534062
$ws.Range("R3").Value = 0.9355806557999999
$ws.Range("S3").Value = 0.3656254573230189
$ws.Range("T3").Value = 0.365625457323019

# --- Row 4: sending cluster FAPs -> MuSCs, target cluster ECs -> FAPs ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 35.755375
$ws.Range("H4").Value = 107.266125
$ws.Range("I4").Value = 0.2275818516336261
$ws.Range("J4").Value = 0.2275818516336262
$ws.Range("M4").Value = 0.001809666666666667
$ws.Range("N4").Value = 0.005429
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.06470531029166667
$ws.Range("R4").Value = 0.582347792625
$ws.Range("S4").Value = 0.2275818516336261
$ws.Range("T4").Value = 0.2275818516336262

# --- Remove the now-obsolete rows 5-7 (MuSCs->ECs, MuSCs->FAPs pairs no longer present) ---
$ws.Rows("5:7").Delete()
